$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update policy number (NroPoliza) in row 2 - keep as text (quote-prefixed string)
$ws.Range("E2").Formula = "'12112002368"

# Update claim date (FechaSiniestro) in row 2 - keep as text (quote-prefixed string)
$ws.Range("G2").Formula = "'19/05/2021"

# Update the active selection to G3 (as recorded in the saved view state)
$ws.Range("G3").Select()
